$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34: copy formatting from row 33 (Medium -> F style matches s="3")
$ws.Range("A33:G33").Copy()
$ws.Range("A34:G34").PasteSpecial(-4122)
$ws.Range("E34").Clear()
$ws.Range("A34").Value = 173
$ws.Range("B34").Value = "Binary Search Tree Iterator"
$ws.Range("C34").Value = "Tree"
$ws.Range("D34").Value = "Aton"
$ws.Range("F34").Value = "Medium"
$ws.Range("G34").Value = "Python"

# Row 35: copy formatting from row 31 (Easy -> F style matches s="6")
$ws.Range("A31:G31").Copy()
$ws.Range("A35:G35").PasteSpecial(-4122)
$ws.Range("E35").Clear()
$ws.Range("A35").Value = 111
$ws.Range("B35").Value = "Minimum Depth of Binary Tree"
$ws.Range("C35").Value = "Tree"
$ws.Range("D35").Value = "Aton"
$ws.Range("F35").Value = "Easy"
$ws.Range("G35").Value = "Python"

# Row 36: copy formatting from row 31 (Easy -> F style matches s="6")
$ws.Range("A31:G31").Copy()
$ws.Range("A36:G36").PasteSpecial(-4122)
$ws.Range("E36").Clear()
$ws.Range("A36").Value = 257
$ws.Range("B36").Value = "Binary Tree Path"
$ws.Range("C36").Value = "Tree"
$ws.Range("D36").Value = "Aton"
$ws.Range("F36").Value = "Easy"
$ws.Range("G36").Value = "Python"

$ws.Range("D33").Select()
$excel.CutCopyMode = 0
